$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.603.87'
$ws.Range("E2").Value = '  -6.84%  '

# Row 3
$ws.Range("D3").Value = '2.960.15'
$ws.Range("E3").Value = '  -8.61%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.32%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.65'
$ws.Range("E5").Value = '  -11.03%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.51'
$ws.Range("E6").Value = '  -15.81%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").Value = '2.941.99'
$ws.Range("E8").Value = '  -9.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.450'
$ws.Range("E9").Value = '  -17.19%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").Value = '  -16.84%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("E11").Value = '  -8.86%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.429'
$ws.Range("E12").Value = '  -12.73%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000203'
$ws.Range("E13").Value = '  -16.51%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.16'
$ws.Range("E14").Value = '  -20.03%  '

# Row 15
$ws.Range("D15").Value = '3.462.20'
$ws.Range("E15").Value = '  -8.14%  '

# Row 16
$ws.Range("D16").Value = '62.484.72'
$ws.Range("E16").Value = '  -7.01%  '

# Row 17
$ws.Range("E17").Value = '  -4.28%  '

# Row 18
$ws.Range("D18").Value = '2.955.29'
$ws.Range("E18").Value = '  -8.79%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '470.64'
$ws.Range("E19").Value = '  -11.93%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.05'
$ws.Range("E20").Value = '  -14.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.39'
$ws.Range("E21").Value = '  -16.99%  '

# Row 22
$ws.Range("E22").Value = '  -16.63%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.38'
$ws.Range("E23").Value = '  -19.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.51'
$ws.Range("E24").Value = '  -13.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.72'
$ws.Range("E25").Value = '  -14.89%  '

# Row 26
$ws.Range("E26").Value = '  -0.22%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.62'
$ws.Range("E27").Value = '  -17.94%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.95'
$ws.Range("E28").Value = '  -14.41%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.83'
$ws.Range("E29").Value = '  -16.39%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '24.35'
$ws.Range("E30").Value = '  -16.52%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.08'
$ws.Range("E31").Value = '  -6.51%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.15%  '

# Row 33
$ws.Range("E33").Value = '  -14.78%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.90'
$ws.Range("E34").Value = '  -3.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '466.70'
$ws.Range("E35").Value = '  -12.96%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.38'
$ws.Range("E36").Value = '  -16.11%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.67'
$ws.Range("E37").Value = '  -18.89%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0382'
$ws.Range("E38").Value = '  -10.21%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0732'
$ws.Range("E39").Value = '  -14.31%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.111'
$ws.Range("E40").Value = '  -9.79%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.64'
$ws.Range("E41").Value = '  -17.62%  '

# Row 42
$ws.Range("D42").Value = '2.622.57'
$ws.Range("E42").Value = '  -10.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("E44").Value = '  -17.99%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.219'
$ws.Range("E45").Value = '  -16.93%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '112.37'
$ws.Range("E46").Value = '  -5.02%  '

# Row 47
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.82'
$ws.Range("E47").Value = '  -16.13%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0988'
$ws.Range("E48").Value = '  -13.91%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.50'
$ws.Range("E49").Value = '  -18.58%  '

# Row 50
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0459'
$ws.Range("E50").Value = '  -22.04%  '

# Row 51
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.18'
$ws.Range("E51").Value = '  -5.43%  '
